$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 28, shifting existing rows 28-118 down to 29-119.
$ws.Rows("28:28").Insert()

# Populate the newly inserted row 28 with the new weekly record.
$ws.Range("A28").Value = 8
$ws.Range("B28").Value = "Terminal La Palmera de La Serena"
$ws.Range("C28").Value = "Coquimbo"
$ws.Range("D28").Value = 44600
$ws.Range("E28").Value = 4
$ws.Range("F28").Value = 100112001
$ws.Range("G28").Value = "Berenjena"
$ws.Range("H28").Value = "Sin especificar"
$ws.Range("I28").Value = "Primera"
$ws.Range("J28").Value = 520
$ws.Range("K28").Value = 8500
$ws.Range("L28").Value = 9000
$ws.Range("M28").Value = 8750
$ws.Range("N28").Value = "$/caja 50 unidades"
$ws.Range("O28").Value = "Región de Arica y Parinacota"
$ws.Range("P28").Value = 175
$ws.Range("Q28").Value = 50
$ws.Range("R28").Value = "Hortaliza"

# Match the date cell style (custom date format) used by the rest of column D.
$ws.Range("D28").NumberFormat = $ws.Range("D29").NumberFormat
